# Commit: "update code bo phan phu cap tai cac co so khac"
# The "Phụ cấp" (allowance) line item is removed from the per-location
# salary breakdown for every location EXCEPT the home office (CẦN THƠ):
#   - row "Phụ cấp tại LONG XUYÊN" is removed
#   - row "Phụ cấp tại SÓC TRĂNG" is removed
# Every row below each deleted row shifts up by one, keeping its own
# label/value pair intact (this is a plain "delete entire row" edit, not a
# content rewrite of the surrounding rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lương")

# Delete the "Phụ cấp tại LONG XUYÊN" row (row 13).
$ws.Rows(13).Delete()

# After the shift above, "Phụ cấp tại SÓC TRĂNG" (originally row 24) is now
# row 23. Delete it too.
$ws.Rows(23).Delete()
